$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two records (rows 9 and 10) got their Id / Antal / Ost / Nord
# values swapped with each other; every other column stays the same.
#   Row 9:  A 111675585 -> 111675587, I "1" -> "3",
#           Q 690349.9096738817 -> 690344.8588249951,
#           R 6661440.004307052 -> 6661440.743740954
#   Row 10: A 111675587 -> 111675585, I "3" -> "1",
#           Q 690344.8588249951 -> 690349.9096738817,
#           R 6661440.743740954 -> 6661440.004307052

# Read the current values first (note: use .Value2 / .Text to actually
# read the contents - the plain .Value getter is not usable here).
$a9 = $ws.Range("A9").Value2
$i9 = $ws.Range("I9").Text
$q9 = $ws.Range("Q9").Value2
$r9 = $ws.Range("R9").Value2

$a10 = $ws.Range("A10").Value2
$i10 = $ws.Range("I10").Text
$q10 = $ws.Range("Q10").Value2
$r10 = $ws.Range("R10").Value2

# Write row 9 with row 10's values (I keeps its original text data type
# by forcing it with a leading apostrophe, same as Excel does).
$ws.Range("A9").Value2 = $a10
$ws.Range("I9").Value2 = "'" + $i10
$ws.Range("Q9").Value2 = $q10
$ws.Range("R9").Value2 = $r10

# Write row 10 with row 9's original values.
$ws.Range("A10").Value2 = $a9
$ws.Range("I10").Value2 = "'" + $i9
$ws.Range("Q10").Value2 = $q9
$ws.Range("R10").Value2 = $r9
